$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changed from 45188 to 45189 for every data row (2..485)
for ($r = 2; $r -le 485; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45189
}

# Row 2 specific updates: a new signal species ("Blomkålssvamp") was found
$ws.Range("I2").Value2 = 1
$ws.Range("Q2").Value2 = 6
$ws.Range("R2").Value2 = "Slåttergubbe`r`nSvinrot`r`nBlomkålssvamp`r`nVanlig padda`r`nMattlummer`r`nRevlummer"

# Keep the row height as it was (the sheet uses a fixed 15pt row height
# everywhere, even for multi-line wrapped text) — writing the longer
# string otherwise triggers autofit in this engine.
$ws.Rows.Item(2).RowHeight = 15
